# Tracker de resultados - actualización automática
# 1) Rellena resultado/profit de apuestas que ya se resolvieron (Fallo / -1)
# 2) Añade la fila nueva con el siguiente partido pendiente de resultado

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Completar resultados pendientes ---
$ws.Range("G133").Value = "Fallo"
$ws.Range("H133").Value = -1

$ws.Range("G134").Value = "Fallo"
$ws.Range("H134").Value = -1

$ws.Range("G136").Value = "Fallo"
$ws.Range("H136").Value = -1

# --- Añadir nueva fila con el próximo evento ---
$newRow = 137
$ws.Range("A" + $newRow).Value = 14807186

# La fecha se guarda como texto ("2025-10-10"), igual que el resto de la
# columna: se antepone un apóstrofo para evitar que Excel la autoconvierta
# a un valor de fecha serial, y luego se restablece el estilo por defecto
# (evita dejar marcado el "quote prefix" en la celda).
$ws.Range("B" + $newRow).Value = "'2025-10-10"
$ws.Range("B" + $newRow).Style = "Normal"

$ws.Range("C" + $newRow).Value = "Jay Dylan Hara Friend"
$ws.Range("D" + $newRow).Value = "Daniel Milavsky"
$ws.Range("E" + $newRow).Value = "Gana Jay Dylan Hara Friend"
$ws.Range("F" + $newRow).Value = 1.57
